$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("per_game")
$ws.Range("G44").Value = 69
$ws.Range("H44").Value = 69
$ws.Range("L44").Value = 0.423
$ws.Range("O44").Value = 0.387
$ws.Range("R44").Value = 0.481
$ws.Range("S44").Value = 0.542
$ws.Range("AE44").Value = 3.4
$ws.Range("G55").Value = 69
$ws.Range("H55").Value = 69
$ws.Range("L55").Value = 0.423
$ws.Range("O55").Value = 0.387
$ws.Range("R55").Value = 0.481
$ws.Range("S55").Value = 0.542
$ws.Range("AE55").Value = 3.4
$ws.Range("G62").Value = 846
$ws.Range("H62").Value = 650
$ws.Range("O62").Value = 0.365
$ws.Range("R62").Value = 0.48
$ws.Range("T62").Value = 0.8
$ws.Range("O64").Value = 0.04199999999999998
$ws.Range("R64").Value = 0.01300000000000001
$ws.Range("T64").Value = -0.3
$ws = $wb.Worksheets.Item("per_minute")
$ws.Range("G34").Value = 69
$ws.Range("H34").Value = 69
$ws.Range("I34").Value = 1779
$ws.Range("K34").Value = 4.1
$ws.Range("L34").Value = 0.423
$ws.Range("O34").Value = 0.387
$ws.Range("P34").Value = 0.7
$ws.Range("R34").Value = 0.481
$ws.Range("G45").Value = 69
$ws.Range("H45").Value = 69
$ws.Range("I45").Value = 1779
$ws.Range("K45").Value = 4.1
$ws.Range("L45").Value = 0.423
$ws.Range("O45").Value = 0.387
$ws.Range("P45").Value = 0.7
$ws.Range("R45").Value = 0.481
$ws.Range("G52").Value = 846
$ws.Range("H52").Value = 650
$ws.Range("I52").Value = 24308
$ws.Range("O52").Value = 0.365
$ws.Range("R52").Value = 0.48
$ws.Range("AD52").Value = 8.5
$ws.Range("O54").Value = 0.04199999999999998
$ws.Range("R54").Value = 0.01300000000000001
$ws.Range("AD54").Value = -0.1999999999999993
$ws = $wb.Worksheets.Item("per_poss")
$ws.Range("G34").Value = 69
$ws.Range("H34").Value = 69
$ws.Range("I34").Value = 1779
$ws.Range("L34").Value = 0.423
$ws.Range("N34").Value = 3.5
$ws.Range("O34").Value = 0.387
$ws.Range("Q34").Value = 2.1
$ws.Range("R34").Value = 0.481
$ws.Range("T34").Value = 0.6
$ws.Range("AC34").Value = 4.7
$ws.Range("AD34").Value = 6.6
$ws.Range("AF34").Value = 117
$ws.Range("AG34").Value = 116
$ws.Range("G45").Value = 69
$ws.Range("H45").Value = 69
$ws.Range("I45").Value = 1779
$ws.Range("L45").Value = 0.423
$ws.Range("N45").Value = 3.5
$ws.Range("O45").Value = 0.387
$ws.Range("Q45").Value = 2.1
$ws.Range("R45").Value = 0.481
$ws.Range("T45").Value = 0.6
$ws.Range("AC45").Value = 4.7
$ws.Range("AD45").Value = 6.6
$ws.Range("AF45").Value = 117
$ws.Range("AG45").Value = 116
$ws.Range("G52").Value = 846
$ws.Range("H52").Value = 650
$ws.Range("I52").Value = 24308
$ws.Range("O52").Value = 0.365
$ws.Range("R52").Value = 0.48
$ws.Range("S52").Value = 1.4
$ws.Range("O54").Value = 0.04199999999999998
$ws.Range("R54").Value = 0.01300000000000001
$ws.Range("S54").Value = -0.4999999999999999
$ws = $wb.Worksheets.Item("advanced")
$ws.Range("G34").Value = 69
$ws.Range("H34").Value = 1779
$ws.Range("I34").Value = 5.7
$ws.Range("J34").Value = 0.5610000000000001
$ws.Range("K34").Value = 0.617
$ws.Range("L34").Value = 0.114
$ws.Range("M34").Value = 6
$ws.Range("N34").Value = 12.2
$ws.Range("S34").Value = 16.3
$ws.Range("X34").Value = 2.4
$ws.Range("Y34").Value = 0.066
$ws.Range("AA34").Value = -3.3
$ws.Range("AB34").Value = 0.3
$ws.Range("AC34").Value = -3
$ws.Range("G45").Value = 69
$ws.Range("H45").Value = 1779
$ws.Range("I45").Value = 5.7
$ws.Range("J45").Value = 0.5610000000000001
$ws.Range("K45").Value = 0.617
$ws.Range("L45").Value = 0.114
$ws.Range("M45").Value = 6
$ws.Range("N45").Value = 12.2
$ws.Range("S45").Value = 16.3
$ws.Range("X45").Value = 2.4
$ws.Range("Y45").Value = 0.066
$ws.Range("AA45").Value = -3.3
$ws.Range("AB45").Value = 0.3
$ws.Range("AC45").Value = -3
$ws.Range("G52").Value = 846
$ws.Range("H52").Value = 24308
$ws.Range("AD52").Value = 6.7
$ws.Range("AD54").Value = -4.7
